$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every B:E cell in this sheet is stored as text in the source file -- including
# cells that look like plain numbers (e.g. "25.00", "19.90", "557.27"). Setting
# .Value directly would let Excel auto-coerce those into numeric cells and drop
# the significant trailing zeros, so each touched cell is forced to the "@"
# (text) number format first. The style is reset to "Normal" afterwards so the
# cell's style index matches the untouched original (only the cached
# NumberFormat differs while the value assignment runs).
$updates = [ordered]@{
    'D2' = '60.398.95'
    'E2' = '  +4.18%  '
    'D3' = '2.451.96'
    'E3' = '  +4.17%  '
    'E4' = '  -0.06%  '
    'D5' = '557.27'
    'E5' = '  +3.17%  '
    'D6' = '139.53'
    'E6' = '  +2.64%  '
    'E7' = '  -0.11%  '
    'D8' = '0.574'
    'E8' = '  +1.60%  '
    'E9' = '  +5.04%  '
    'D10' = '5.83'
    'E10' = '  +4.55%  '
    'D11' = '0.363'
    'E11' = '  +2.64%  '
    'E12' = '  -1.93%  '
    'D13' = '25.02'
    'E13' = '  +4.88%  '
    'D14' = '2.888.13'
    'E14' = '  +4.14%  '
    'D15' = '60.299.71'
    'E15' = '  +4.06%  '
    'E16' = '  +5.45%  '
    'D17' = '2.449.74'
    'E17' = '  +3.84%  '
    'D18' = '11.52'
    'E18' = '  +7.67%  '
    'E19' = '  +4.20%  '
    'D20' = '336.52'
    'E20' = '  +1.33%  '
    'E21' = '  +2.40%  '
    'D23' = '64.88'
    'E23' = '  +3.08%  '
    'E24' = '  +2.32%  '
    'D25' = '8.63'
    'E25' = '  +1.54%  '
    'E26' = '  +0.05%  '
    'E27' = '  +0.65%  '
    'D28' = '0.0₃0800'
    'E28' = '  +8.67%  '
    'D29' = '1.82'
    'E29' = '  +3.89%  '
    'E30' = '  +3.04%  '
    'D31' = '170.77'
    'E31' = '  -0.95%  '
    'D32' = '18.89'
    'E32' = '  +2.28%  '
    'E33' = '  -0.28%  '
    'E35' = '  +6.29%  '
    'D36' = '4.32'
    'E36' = '  +2.37%  '
    'E37' = '  +0.14%  '
    'D38' = '1.65'
    'E38' = '  +0.73%  '
    'D39' = '40.18'
    'E39' = '  +2.20%  '
    'D40' = '0.421'
    'E40' = '  +11.33%  '
    'D41' = '317.67'
    'E41' = '  +7.91%  '
    'B42' = 'Aave'
    'C42' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D42' = '145.27'
    'E42' = '  -0.35%  '
    'B43' = 'Filecoin'
    'C43' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D43' = '3.75'
    'E43' = '  +2.59%  '
    'B44' = 'Hedera'
    'C44' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D44' = '0.0528'
    'E44' = '  +4.86%  '
    'B45' = 'Stellar'
    'C45' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D45' = '0.0965'
    'E45' = '  +1.61%  '
    'D46' = '19.90'
    'E46' = '  +3.18%  '
    'B47' = 'Mantle'
    'C47' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D47' = '0.576'
    'E47' = '  +2.09%  '
    'B48' = 'Polygon'
    'C48' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'D48' = '0.408'
    'E48' = '  +7.05%  '
    'D49' = '0.0228'
    'E49' = '  +3.00%  '
    'D50' = '11.04'
    'E50' = '  -0.15%  '
    'D51' = '1.65'
    'E51' = '  +6.41%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $updates[$cellRef]
    $ws.Range($cellRef).Style = "Normal"
}
